# The "Recorded By" column (G) lists the users who recorded/edited each
# attendance session as a comma-separated string. This edit reverses the
# order of the names in that list for every data row on the active sheet
# (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System").
# Single-value cells are left unchanged (reversing one item is a no-op).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $newVal = $reversed -join ", "
            $cell.Value = $newVal
        }
    }
}
